$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.410.62'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '1.805.03'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'227.65"
$ws.Range("E5").Value = '  +0.53%  '
$ws.Range("E6").Value = '  +4.10%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = "'35.85"
$ws.Range("E8").Value = '  +8.68%  '
$ws.Range("E9").Value = '  +2.50%  '
$ws.Range("D10").Value = "'0.0695"
$ws.Range("E10").Value = '  +0.97%  '
$ws.Range("E11").Value = '  +2.12%  '
$ws.Range("D12").Value = '2.066.33'
$ws.Range("E12").Value = '  +1.06%  '
$ws.Range("D13").Value = "'11.43"
$ws.Range("E13").Value = '  +2.03%  '
$ws.Range("D14").Value = '1.806.69'
$ws.Range("E14").Value = '  +1.21%  '
$ws.Range("E16").Value = '  +5.34%  '
$ws.Range("D17").Value = '34.402.84'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").Value = "'69.25"
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("D20").Value = '0.0₃0798'
$ws.Range("E20").Value = '  +0.39%  '
$ws.Range("E21").Value = '  +1.75%  '
$ws.Range("E23").Value = '  +1.07%  '
$ws.Range("E24").Value = '  +3.26%  '
$ws.Range("D25").Value = "'170.88"
$ws.Range("E25").Value = '  +1.41%  '
$ws.Range("D26").Value = "'7.94"
$ws.Range("E26").Value = '  +8.24%  '
$ws.Range("D27").Value = "'16.97"
$ws.Range("E27").Value = '  +2.70%  '
$ws.Range("E28").Value = '  +3.10%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").Value = "'4.07"
$ws.Range("E30").Value = '  +1.03%  '
$ws.Range("E31").Value = '  +1.09%  '
$ws.Range("E32").Value = '  +1.72%  '
$ws.Range("E33").Value = '  +0.80%  '
$ws.Range("E34").Value = '  +0.99%  '
$ws.Range("D35").Value = '1.398.36'
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("D36").Value = "'0.673"
$ws.Range("E36").Value = '  -1.72%  '
$ws.Range("E37").Value = '  -3.49%  '
$ws.Range("D38").Value = "'1.06"
$ws.Range("E38").Value = '  -0.47%  '
$ws.Range("E39").Value = '  -0.19%  '
$ws.Range("E40").Value = '  +11.49%  '
$ws.Range("D41").Value = "'0.965"
$ws.Range("E41").Value = '  +2.76%  '
$ws.Range("D42").Value = "'82.69"
$ws.Range("E42").Value = '  -2.36%  '
$ws.Range("D43").Value = "'2.83"
$ws.Range("E43").Value = '  +1.99%  '
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("D45").Value = "'13.48"
$ws.Range("E45").Value = '  -3.76%  '
$ws.Range("D46").Value = "'6.04"
$ws.Range("E46").Value = '  -0.70%  '
$ws.Range("D47").Value = "'0.0502"
$ws.Range("E47").Value = '  -4.99%  '
$ws.Range("D48").Value = '1.966.56'
$ws.Range("E48").Value = '  +1.04%  '
$ws.Range("D49").Value = "'105.24"
$ws.Range("E49").Value = '  -0.11%  '
$ws.Range("E51").Value = '  +1.30%  '
